# ProVisioNET_Poster_Lehrpersonen.pptx edits
#
# 1. Title shape (id=2, "Titel 1"): the subtitle line had "Visi" and "o"
#    as two separate runs (first bold, second not). The author retyped
#    across the run boundary so they collapsed into a single bold run
#    "Visio".
# 2. Rectangle shape (id=18, "Rechteck 17"): fixed a typo, "Ihrem" ->
#    "Ihren", inside the single run that holds the whole sentence.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# --- 1. Title shape: merge "Visi" + "o" runs into a single "Visio" run ---
$title = Get-ShapeById $s 2
$titleRange = $title.TextFrame.TextRange
$titleText = $titleRange.Text

$firstVisi = $titleText.IndexOf("Visi")
$secondVisi = $titleText.IndexOf("Visi", $firstVisi + 1)

if ($secondVisi -ge 0) {
    # 1-based character index of the "V" in the second "Visi" occurrence,
    # spanning the 5 characters "Visio" (the run "Visi" plus the
    # following single-character run "o").
    $start = $secondVisi + 1
    $visioRange = $titleRange.Characters($start, 5)
    $visioRange.Text = "Visio"
}

# --- 2. Rectangle shape: fix "Ihrem" -> "Ihren" ---
$rect = Get-ShapeById $s 18
$rectRange = $rect.TextFrame.TextRange
$rectText = $rectRange.Text
$rectRange.Text = $rectText.Replace("Ihrem Blick", "Ihren Blick")

# --- 3. Best-effort: refresh the cached "datetimeFigureOut" text on the
#    handout master / notes master date placeholders (17.01.2022 ->
#    28.01.2022). PowerPoint recomputes this automatically from the
#    system clock whenever the deck is saved; it is not a deliberate
#    text edit, so this is wrapped defensively and must never abort the
#    rest of the script if the host does not support writing to a
#    field-backed TextRange.
try {
    $handoutMaster = $p.HandoutMaster
    for ($i = 1; $i -le $handoutMaster.Shapes.Count; $i++) {
        $phShape = $handoutMaster.Shapes.Item($i)
        if ($phShape.Name -like "Datumsplatzhalter*") {
            $phShape.TextFrame.TextRange.Text = "28.01.2022"
        }
    }
} catch {
}

try {
    $notesMaster = $p.NotesMaster
    for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
        $phShape = $notesMaster.Shapes.Item($i)
        if ($phShape.Name -like "Datumsplatzhalter*") {
            $phShape.TextFrame.TextRange.Text = "28.01.2022"
        }
    }
} catch {
}
